# Apply the crypto price/volume updates for this commit.
# Preserve text semantics of the Price column (D) by forcing a
# text number format before assigning number-looking strings, so
# Excel does not silently coerce values like "2.80" or "1.00"
# into numeric 2.8 / 1, which would lose the original formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "52.053.35"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "2.797.55"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "362.23"
$ws.Range("E5").Value = "  +2.45%  "
$ws.Range("D6").Value = "111.13"
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("D7").Value = "0.563"
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.595"
$ws.Range("E9").Value = "  -1.48%  "
$ws.Range("D10").Value = "40.15"
$ws.Range("E10").Value = "  -2.53%  "
$ws.Range("D11").Value = "0.0858"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("E12").Value = "  +1.62%  "
$ws.Range("D13").Value = "19.51"
$ws.Range("E13").Value = "  -1.64%  "
$ws.Range("D14").Value = "7.65"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").Value = "3.238.77"
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("D16").Value = "2.810.86"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "0.955"
$ws.Range("E17").Value = "  +7.67%  "
$ws.Range("D18").Value = "52.025.61"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").Value = "7.49"
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("D21").Value = "13.18"
$ws.Range("E21").Value = "  -1.70%  "
$ws.Range("D22").Value = "0.0₃0985"
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("D23").Value = "273.32"
$ws.Range("E23").Value = "  +1.75%  "
$ws.Range("E24").Value = "  +1.42%  "
$ws.Range("D25").Value = "2.80"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").Value = "26.70"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.157"
$ws.Range("E28").Value = "  +13.39%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "10.26"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("D31").Value = "52.31"
$ws.Range("E31").Value = "  +1.26%  "
$ws.Range("D32").Value = "0.0466"
$ws.Range("E32").Value = "  +4.20%  "
$ws.Range("D33").Value = "34.36"
$ws.Range("E33").Value = "  +1.96%  "
$ws.Range("D34").Value = "5.81"
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("D35").Value = "0.0855"
$ws.Range("E35").Value = "  +3.71%  "
$ws.Range("D36").Value = "5.27"
$ws.Range("E36").Value = "  +1.39%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "3.23"
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("D39").Value = "18.51"
$ws.Range("E39").Value = "  +1.78%  "
$ws.Range("E40").Value = "  -2.21%  "
$ws.Range("E41").Value = "  +4.19%  "
$ws.Range("D42").Value = "0.116"
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "2.25"
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "121.48"
$ws.Range("E44").Value = "  -3.07%  "
$ws.Range("D45").Value = "22.32"
$ws.Range("E45").Value = "  -4.96%  "
$ws.Range("D46").Value = "2.078.77"
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("D47").Value = "3.26"
$ws.Range("E47").Value = "  -2.20%  "
$ws.Range("D48").Value = "2.21"
$ws.Range("E48").Value = "  -1.89%  "
$ws.Range("D49").Value = "5.74"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").Value = "0.950"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("D51").Value = "8.94"
$ws.Range("E51").Value = "  +0.27%  "
